$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'isophonics_165'
$ws.Cells.Item(2, 2).Value = 'schubert-winterreise_141'
$ws.Cells.Item(2, 3).Value = 0.115546218487395
$ws.Cells.Item(2, 4).Value = '[[''A'', ''E'', ''A'']]'
$ws.Cells.Item(2, 5).Value = '[[''F:maj'', ''C:maj'', ''F:maj'']]'
$ws.Cells.Item(2, 6).Value = '[(52.680839, 55.420793)]'
$ws.Cells.Item(2, 7).Value = '[(0.78, 24.86)]'
$ws.Cells.Item(2, 8).Value = ''

# Row 3
$ws.Cells.Item(3, 1).Value = 'isophonics_152'
$ws.Cells.Item(3, 2).Value = 'schubert-winterreise_74'
$ws.Cells.Item(3, 3).Value = 0.1576923076923077
$ws.Cells.Item(3, 4).Value = '[[''A/3'', ''D'', ''A'']]'
$ws.Cells.Item(3, 5).Value = '[[''F:maj'', ''A#:maj'', ''F:maj'']]'
$ws.Cells.Item(3, 6).Value = '[(4.25, 9.179)]'
$ws.Cells.Item(3, 7).Value = '[(129.38, 134.38)]'
$ws.Cells.Item(3, 9).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

# Row 4
$ws.Cells.Item(4, 1).Value = 'isophonics_74'
$ws.Cells.Item(4, 2).Value = 'isophonics_50'
$ws.Cells.Item(4, 3).Value = 0.275
$ws.Cells.Item(4, 4).Value = '[[''A'', ''E'', ''A''], [''E'', ''A'', ''E'']]'
$ws.Cells.Item(4, 5).Value = '[[''C'', ''G'', ''C''], [''G'', ''C'', ''G'']]'
$ws.Cells.Item(4, 6).Value = '[(35.197913, 43.754467), (37.995918, 46.610521)]'
$ws.Cells.Item(4, 7).Value = '[(58.300068, 63.443287), (37.936167, 42.951678)]'

# Row 5
$ws.Cells.Item(5, 1).Value = 'schubert-winterreise_141'
$ws.Cells.Item(5, 2).Value = 'isophonics_200'
$ws.Cells.Item(5, 3).Value = 0.1366459627329192
$ws.Cells.Item(5, 4).Value = '[[''F:maj'', ''C:maj'', ''F:maj'']]'
$ws.Cells.Item(5, 5).Value = '[[''A'', ''E'', ''A'']]'
$ws.Cells.Item(5, 6).Value = '[(0.78, 24.86)]'
$ws.Cells.Item(5, 7).Value = '[(10.953139, 19.672232)]'

# Row 6
$ws.Cells.Item(6, 1).Value = 'jaah_65'
$ws.Cells.Item(6, 2).Value = 'jaah_19'
$ws.Cells.Item(6, 3).Value = 0.02827442827442828
$ws.Cells.Item(6, 4).Value = '[[''F:7'', ''Bb:min7'', ''Eb:7'', ''Ab:maj6''], [''Db:min7'', ''Gb:7'', ''B:maj7'', ''E:min7'']]'
$ws.Cells.Item(6, 5).Value = '[[''F:7'', ''Bb:min7'', ''Eb:7'', ''Ab:maj6''], [''C#:min7'', ''F#:7'', ''B:maj7'', ''E:min7'']]'
$ws.Cells.Item(6, 6).Value = '[(7.42, 9.88), (20.26, 23.37)]'
$ws.Cells.Item(6, 7).Value = '[(32.74, 35.48), (17.85, 20.65)]'

# Row 7
$ws.Cells.Item(7, 1).Value = 'schubert-winterreise_215'
$ws.Cells.Item(7, 2).Value = 'schubert-winterreise_169'
$ws.Cells.Item(7, 3).Value = 0.1517241379310345
$ws.Cells.Item(7, 4).Value = '[[''G:min'', ''D:7/G'', ''G:min'']]'
$ws.Cells.Item(7, 5).Value = '[[''G:min'', ''D:7'', ''G:min'']]'
$ws.Cells.Item(7, 6).Value = '[(15.78, 21.28)]'
$ws.Cells.Item(7, 7).Value = '[(19.44, 28.3)]'
$ws.Cells.Item(7, 9).Value = 'spotify:track:3OD2uwEUQKg0WyW9Lewata'

# Row 8
$ws.Cells.Item(8, 1).Value = 'schubert-winterreise_48'
$ws.Cells.Item(8, 2).Value = 'schubert-winterreise_44'
$ws.Cells.Item(8, 3).Value = 0.09642857142857142
$ws.Cells.Item(8, 4).Value = '[[''F:maj/C'', ''C:7'', ''F:maj'']]'
$ws.Cells.Item(8, 5).Value = '[[''A#/F'', ''F:7'', ''A#'']]'
$ws.Cells.Item(8, 6).Value = '[(59.5, 65.04)]'
$ws.Cells.Item(8, 7).Value = '[(271.22, 275.78)]'

# Row 9
$ws.Cells.Item(9, 1).Value = 'isophonics_193'
$ws.Cells.Item(9, 2).Value = 'schubert-winterreise_82'
$ws.Cells.Item(9, 3).Value = 0.1440993788819876
$ws.Cells.Item(9, 4).Value = '[[''Bb:7'', ''Eb'', ''Bb'', ''Eb'']]'
$ws.Cells.Item(9, 5).Value = '[[''D:7'', ''G:maj'', ''D:maj'', ''G:maj'']]'
$ws.Cells.Item(9, 6).Value = '[(18.253, 22.735)]'
$ws.Cells.Item(9, 7).Value = '[(44.16, 47.28)]'
$ws.Cells.Item(9, 8).Value = ''
$ws.Cells.Item(9, 9).Value = 'spotify:track:4lrfYSnZmpXdCWuWqVo8L0'

# Row 10
$ws.Cells.Item(10, 1).Value = 'schubert-winterreise_41'
$ws.Cells.Item(10, 2).Value = 'schubert-winterreise_94'
$ws.Cells.Item(10, 3).Value = 0.2913752913752914
$ws.Cells.Item(10, 4).Value = '[[''B:min'', ''F#:maj/A#'', ''B:min'', ''B:7/A'', ''E:min/G'']]'
$ws.Cells.Item(10, 5).Value = '[[''A#:min'', ''F:maj'', ''A#:min'', ''A#:7'', ''D#:min/A#'']]'
$ws.Cells.Item(10, 6).Value = '[(13.74, 17.58)]'
$ws.Cells.Item(10, 7).Value = '[(15.6, 27.52)]'
$ws.Cells.Item(10, 8).Value = ''
$ws.Cells.Item(10, 9).Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

# Row 11
$ws.Cells.Item(11, 1).Value = 'schubert-winterreise_132'
$ws.Cells.Item(11, 2).Value = 'isophonics_48'
$ws.Cells.Item(11, 3).Value = 0.1113122171945701
$ws.Cells.Item(11, 4).Value = '[[''F:maj'', ''A#:maj/F'', ''F:maj'']]'
$ws.Cells.Item(11, 5).Value = '[[''A'', ''D'', ''A'']]'
$ws.Cells.Item(11, 6).Value = '[(17.34, 24.4)]'
$ws.Cells.Item(11, 7).Value = '[(82.290863, 90.185647)]'
$ws.Cells.Item(11, 8).Value = 'spotify:track:68YORkKP9uvlOQFMZZZwH5'

# Row 12
$ws.Cells.Item(12, 1).Value = 'schubert-winterreise_109'
$ws.Cells.Item(12, 2).Value = 'schubert-winterreise_146'
$ws.Cells.Item(12, 3).Value = 0.09779367918902802
$ws.Cells.Item(12, 4).Value = '[[''F:maj/C'', ''C:7'', ''F:maj'', ''D:min/A'']]'
$ws.Cells.Item(12, 5).Value = '[[''D:maj/F#'', ''A:7'', ''D:maj'', ''B:min'']]'
$ws.Cells.Item(12, 6).Value = '[(102.16, 113.9)]'
$ws.Cells.Item(12, 7).Value = '[(38.62, 46.84)]'
$ws.Cells.Item(12, 8).Value = 'spotify:track:5UYEp9kllA47IhttiiMuJ0'
$ws.Cells.Item(12, 9).Value = 'spotify:track:2g41AZ58LFdQLxmWx82ujI'

# Row 13
$ws.Cells.Item(13, 1).Value = 'isophonics_265'
$ws.Cells.Item(13, 2).Value = 'schubert-winterreise_129'
$ws.Cells.Item(13, 3).Value = 0.2871794871794872
$ws.Cells.Item(13, 4).Value = '[[''D'', ''A'', ''D'', ''A'']]'
$ws.Cells.Item(13, 5).Value = '[[''A#:maj'', ''F:maj'', ''A#:maj'', ''F:maj'']]'
$ws.Cells.Item(13, 6).Value = '[(102.174263, 115.003287)]'
$ws.Cells.Item(13, 7).Value = '[(107.64, 108.84)]'
$ws.Cells.Item(13, 9).Value = 'spotify:track:1nvxQGWCnikMK7a4HYQvSx'

# Row 14
$ws.Cells.Item(14, 1).Value = 'isophonics_21'
$ws.Cells.Item(14, 2).Value = 'isophonics_28'
$ws.Cells.Item(14, 3).Value = 0.3091787439613526
$ws.Cells.Item(14, 4).Value = '[[''G/3'', ''C'', ''G/3'', ''C'']]'
$ws.Cells.Item(14, 5).Value = '[[''G'', ''C/5'', ''G'', ''C'']]'
$ws.Cells.Item(14, 6).Value = '[(130.528, 132.075)]'
$ws.Cells.Item(14, 7).Value = '[(64.476575, 70.099405)]'
$ws.Cells.Item(14, 9).Value = ''

# Row 15
$ws.Cells.Item(15, 1).Value = 'isophonics_159'
$ws.Cells.Item(15, 2).Value = 'isophonics_30'
$ws.Cells.Item(15, 3).Value = 0.1332236842105263
$ws.Cells.Item(15, 4).Value = '[[''A'', ''G/2'', ''D/5'']]'
$ws.Cells.Item(15, 5).Value = '[[''A'', ''G'', ''D'']]'
$ws.Cells.Item(15, 6).Value = '[(59.723, 63.111)]'
$ws.Cells.Item(15, 7).Value = '[(7.645256, 11.062412)]'

# Row 16
$ws.Cells.Item(16, 1).Value = 'isophonics_19'
$ws.Cells.Item(16, 2).Value = 'isophonics_99'
$ws.Cells.Item(16, 3).Value = 0.1375
$ws.Cells.Item(16, 4).Value = '[[''G:min'', ''C'', ''F'']]'
$ws.Cells.Item(16, 5).Value = '[[''F#:min'', ''B'', ''E'']]'
$ws.Cells.Item(16, 6).Value = '[(45.474603, 48.945986)]'
$ws.Cells.Item(16, 7).Value = '[(16.376439, 18.860975)]'
$ws.Cells.Item(16, 8).Value = ''
$ws.Cells.Item(16, 9).Value = ''

# Row 17
$ws.Cells.Item(17, 1).Value = 'isophonics_156'
$ws.Cells.Item(17, 2).Value = 'isophonics_216'
$ws.Cells.Item(17, 3).Value = 0.1821631878557875
$ws.Cells.Item(17, 4).Value = '[[''E'', ''A'', ''D'', ''A'']]'
$ws.Cells.Item(17, 5).Value = '[[''E'', ''A'', ''D'', ''A'']]'
$ws.Cells.Item(17, 6).Value = '[(9.30541, 13.723731)]'
$ws.Cells.Item(17, 7).Value = '[(25.850181, 33.826235)]'
